# Automatische test-sync: 2025-08-19 19:54:50
# Append a new log row (row 13) to the "Logs" sheet, extend the
# conditional-formatting ranges to cover it, and bump the matching
# "Dashboard" summary count.

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A13").Value = "Vraag over product"
$logs.Range("B13").Value = "documentatie@testbedrijf123.nl"
$logs.Range("D13").Value = "Intern verzoek / Actie voor medewerker"
$logs.Range("F13").Value = "2025-08-19 19:54:27"
$logs.Range("G13").Value = "Nee"
$logs.Range("H13").Value = "Ja"
$logs.Range("I13").Value = "Nee"
$logs.Range("J13").Value = "Nee"

# Extend the existing conditional formatting (D/G/H/I/J, rows 2-12) so it
# also covers the newly added row 13.
$cols = @("D", "G", "H", "I", "J")
foreach ($col in $cols) {
    $existing = $logs.Range($col + "2")
    $fcs = $existing.FormatConditions
    $newRange = $logs.Range($col + "2:" + $col + "13")
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fc = $fcs.Item($i)
        $fc.ModifyAppliesToRange($newRange)
    }
}

# Update the Dashboard summary count for this category.
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Range("B2").Value = 12
